# Fix the mislabeled 2050 column header (it was accidentally left holding a
# stray numeric value instead of the "2050" / "2041-2050" text label used by
# its neighbouring header cells) and drop the "Total" summary rows that were
# appended to each table.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    # Positional params - named binding ($cell, $text) is unreliable here.
    param($cell, $text)
    # Writing a plain string via .Value lets Excel "smart type" a
    # numeric-looking string (e.g. "2050") back into a number. Routing the
    # literal through a formula that evaluates to a string forces text type
    # regardless of the cell's number format, then we collapse the formula
    # back down to a plain value in place - this preserves the cell's
    # existing style (border/bold/center) without minting a new one.
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163) # xlPasteValues
}

# --- Sheet 1: Potencia Acumulada - SIN (MW) ---
$ws = $wb.Worksheets.Item("Potencia Acumulada - SIN (MW)")
Set-TextValue $ws.Range("E1") "2050"
$ws.Rows.Item(13).Delete()

# --- Sheet 2: Geracao Periodo Medio (MWMed) ---
$ws = $wb.Worksheets.Item("Geracao Periodo Medio (MWMed)")
Set-TextValue $ws.Range("E1") "2050"
$ws.Rows.Item(13).Delete()

# --- Sheet 3: Atendimento a Ponta(MW) ---
$ws = $wb.Worksheets.Item("Atendimento a Ponta(MW)")
Set-TextValue $ws.Range("E1") "2050"
$ws.Rows.Item(13).Delete()

# --- Sheet 4: Potencia Incremental - SIN(MW) ---
$ws = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
Set-TextValue $ws.Range("E1") "2041-2050"
$ws.Rows.Item(13).Delete()

# --- Sheet 5: Emissoes Totais (MtCO2eq) ---
$ws = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")
Set-TextValue $ws.Range("E1") "2050"

# --- Sheet 6: Custo Total (bilhões de R$) ---
$ws = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws.Rows.Item(4).Delete()
